$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: remove the thin bottom border (table no longer ends here) ---
$ws.Range("A4:E4").Borders.Item(9).LineStyle = -4142

# Taller row now that A4 carries wrapped text
$ws.Rows.Item(4).RowHeight = 43.2

# A4 gets a new script filename entry
$ws.Range("A4").Value = "SCRIPT/P02P01A/um1102.ssb"

# --- Row 5: new script filename row, formatted like A2/A3 ---
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Rows.Item(5).RowHeight = 43.2
$ws.Range("A5").Value = "SCRIPT/P02P01A/um1105.ssb"

# --- Row 6: new script filename row, formatted like A2/A3 ---
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Range("A6").Value = "SCRIPT/P02P01A/um1108.ssb"

$excel.CutCopyMode = $false

# Move the selection like the saved workbook shows
[void]$ws.Range("D4").Select()
